# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" fund-holdings sheet right after "2021-Q4"
#   (reusing the sheetId that "总计" used to have).
# - Recreate the "总计" (totals) sheet after it, with a new row for
#   2022-Q1 on top and all the other rows shifted down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Reference sheet used purely as a style/format donor (border + bold +
# centered header style, and the "index column" style) so the new
# sheets pick up the very same style ids already used across the
# workbook instead of creating brand-new ones.
# ---------------------------------------------------------------------
$styleDonor = $wb.Worksheets.Item("2021-Q4")

# Grab the old "总计" sheet before it disappears - we need its data to
# rebuild the totals table in the new sheet.
$oldTotal = $wb.Worksheets.Item("总计")

$oldDates  = @()
$oldCounts = @()
$oldValues = @()
for ($r = 2; $r -le 6; $r++) {
    $oldDates  += $oldTotal.Cells.Item($r, 2).Value()
    $oldCounts += $oldTotal.Cells.Item($r, 3).Value()
    $oldValues += $oldTotal.Cells.Item($r, 4).Value()
}

# Remove the old totals sheet - this frees up sheetId 6 so the next
# new sheet we create reuses it (matching sheetId="6" in the target).
$oldTotal.Delete()

# ---------------------------------------------------------------------
# 1) New "2022-Q1" sheet, inserted right after "2021-Q4"
# ---------------------------------------------------------------------
$q1 = $wb.Worksheets.Add($null, $styleDonor)
$q1.Name = "2022-Q1"

$q1.PageSetup.LeftMargin   = 0.75 * 72
$q1.PageSetup.RightMargin  = 0.75 * 72
$q1.PageSetup.TopMargin    = 1 * 72
$q1.PageSetup.BottomMargin = 1 * 72
$q1.PageSetup.HeaderMargin = 0.5 * 72
$q1.PageSetup.FooterMargin = 0.5 * 72

# Copy the header-row + index-column formatting from the donor sheet.
$styleDonor.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)
$styleDonor.Range("A2").Copy()
$q1.Range("A2:A4").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Index column (0-based row counter), mirrors the other quarter sheets.
$q1.Range("A2").Value = 0
$q1.Range("A3").Value = 1
$q1.Range("A4").Value = 2

# Column B (fund code) and columns D:G hold numeric-looking text (e.g.
# "590008", "9.00") stored as text, exactly like the other quarter
# sheets - force text formatting first so Excel doesn't silently
# coerce them into numbers (and lose leading zeros on fund codes).
$q1.Range("B2:B4").NumberFormat = "@"
$q1.Range("D2:G4").NumberFormat = "@"

$q1.Range("B2").Value = "590008"
$q1.Range("C2").Value = "中邮战略新兴产业混合"
$q1.Range("D2").Value = "9.00"
$q1.Range("E2").Value = "88.88"
$q1.Range("F2").Value = "3.35"
$q1.Range("G2").Value = "0.3015"
$q1.Range("H2").Value = 8

$q1.Range("B3").Value = "001571"
$q1.Range("C3").Value = "嘉合磐石混合A"
$q1.Range("D3").Value = "0.64"
$q1.Range("E3").Value = "39.83"
$q1.Range("F3").Value = "3.37"
$q1.Range("G3").Value = "0.0216"
$q1.Range("H3").Value = 7

$q1.Range("B4").Value = "001572"
$q1.Range("C4").Value = "嘉合磐石混合C"
$q1.Range("D4").Value = "0.51"
$q1.Range("E4").Value = "39.83"
$q1.Range("F4").Value = "3.37"
$q1.Range("G4").Value = "0.0172"
$q1.Range("H4").Value = 7

# ---------------------------------------------------------------------
# 2) New "总计" sheet, inserted right after "2022-Q1", rebuilding the
#    totals table with the 2022-Q1 row prepended.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

$total.PageSetup.LeftMargin   = 0.75 * 72
$total.PageSetup.RightMargin  = 0.75 * 72
$total.PageSetup.TopMargin    = 1 * 72
$total.PageSetup.BottomMargin = 1 * 72
$total.PageSetup.HeaderMargin = 0.5 * 72
$total.PageSetup.FooterMargin = 0.5 * 72

$styleDonor.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$styleDonor.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$dates  = @("2022-Q1") + $oldDates
$counts = @(3) + $oldCounts
$values = @(0.34) + $oldValues

for ($i = 0; $i -lt $dates.Length; $i++) {
    $r = $i + 2
    $total.Cells.Item($r, 1).Value = $i
    $total.Cells.Item($r, 2).Value = $dates[$i]
    $total.Cells.Item($r, 3).Value = $counts[$i]
    $total.Cells.Item($r, 4).Value = $values[$i]
}

$wb.Save()
